$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor - update values only
$ws.Range("B3").Value = 0.84457726188958
$ws.Range("C3").Value = 0.8482158566782154
$ws.Range("D3").Value = 0.8025913851327323

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor, update values
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.8445631312607134
$ws.Range("C4").Value = 0.8445631312607134
$ws.Range("D4").Value = 0.7913798731503779

# Row 5: AdaBoostRegressor -> MLPRegressor, update values
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.6533590851872334
$ws.Range("C5").Value = 0.6613432097973514
$ws.Range("D5").Value = 0.620851658905046
